$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row to the table (Table1), extending it to A1:E13
$table = $ws.ListObjects.Item("Table1")
$newRow = $table.ListRows.Add()

$ws.Range("A13").Value = 45431
$ws.Range("B13").Value = "Spez 1"
$ws.Range("C13").Value = "Kaffee und Gipfeli"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 3

# Update selection to match post-edit state
$ws.Range("D18").Select()
